# Auto-generated: update price/profit columns (H:N) per scheduled-runner data refresh
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1410.6154
$ws.Range("I42").Value = 123
$ws.Range("J42").Value = 2514.2856
$ws.Range("K42").Value = 369
$ws.Range("L42").Value = 7542.8568
$ws.Range("M42").Value = -139
$ws.Range("N42").Value = -8002.8568
$ws.Range("H70").Value = 1355.1666
$ws.Range("I70").Value = 1199
$ws.Range("J70").Value = 1433.25
$ws.Range("K70").Value = 3597
$ws.Range("L70").Value = 4299.75
$ws.Range("M70").Value = -3327
$ws.Range("N70").Value = -4839.75
$ws.Range("H73").Value = 1355.1666
$ws.Range("I73").Value = 1199
$ws.Range("J73").Value = 1433.25
$ws.Range("K73").Value = 3597
$ws.Range("L73").Value = 4299.75
$ws.Range("M73").Value = -2661
$ws.Range("N73").Value = -6171.75
$ws.Range("H100").Value = 1970.5333
$ws.Range("I100").Value = 762
$ws.Range("J100").Value = 2574.8
$ws.Range("K100").Value = 762
$ws.Range("L100").Value = 2574.8
$ws.Range("M100").Value = -221
$ws.Range("N100").Value = -3656.8
$ws.Range("H141").Value = 3455.3914
$ws.Range("I141").Value = 4087.111
$ws.Range("K141").Value = 12261.333
$ws.Range("M141").Value = -7081.332999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5655.722
$ws.Range("I74").Value = 2932.1724
$ws.Range("J74").Value = 16939
$ws.Range("K74").Value = 2932.1724
$ws.Range("L74").Value = 16939
$ws.Range("M74").Value = -2058.1724
$ws.Range("N74").Value = -18687
$ws.Range("H77").Value = 5655.722
$ws.Range("I77").Value = 2932.1724
$ws.Range("J77").Value = 16939
$ws.Range("K77").Value = 14660.862
$ws.Range("L77").Value = 84695
$ws.Range("M77").Value = -10292.862
$ws.Range("N77").Value = -93431
$ws.Range("H102").Value = 2878.3333
$ws.Range("I102").Value = 2292.5
$ws.Range("K102").Value = 2292.5
$ws.Range("M102").Value = -670.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1238.6177
$ws.Range("I94").Value = 1164.1428
$ws.Range("J94").Value = 1358.9231
$ws.Range("K94").Value = 1164.1428
$ws.Range("L94").Value = 1358.9231
$ws.Range("M94").Value = -713.1428000000001
$ws.Range("N94").Value = -2260.9231
$ws.Range("H99").Value = 2411.353
$ws.Range("I99").Value = 2562.125
$ws.Range("J99").Value = 2277.3333
$ws.Range("K99").Value = 2562.125
$ws.Range("L99").Value = 2277.3333
$ws.Range("M99").Value = -1064.125
$ws.Range("N99").Value = -5273.3333
$ws.Range("H105").Value = 6870.9688
$ws.Range("I105").Value = 3227.8262
$ws.Range("K105").Value = 3227.8262
$ws.Range("M105").Value = -1480.8262
$ws.Range("H107").Value = 2601.9167
$ws.Range("I107").Value = 2630
$ws.Range("J107").Value = 2562.6
$ws.Range("K107").Value = 2630
$ws.Range("L107").Value = 2562.6
$ws.Range("M107").Value = -710
$ws.Range("N107").Value = -6402.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1438.7646
$ws.Range("I94").Value = 1661.5714
$ws.Range("J94").Value = 1282.8
$ws.Range("K94").Value = 1661.5714
$ws.Range("L94").Value = 1282.8
$ws.Range("M94").Value = -1210.5714
$ws.Range("N94").Value = -2184.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 21322
$ws.Range("J51").Value = 21322
$ws.Range("L51").Value = 63966
$ws.Range("N51").Value = -64886
$ws.Range("H52").Value = 1266.5
$ws.Range("J52").Value = 1266.5
$ws.Range("L52").Value = 3799.5
$ws.Range("N52").Value = -4331.5
$ws.Range("H122").Value = 613.8570999999999
$ws.Range("I122").Value = 445.86667
$ws.Range("J122").Value = 1033.8334
$ws.Range("K122").Value = 4012.80003
$ws.Range("L122").Value = 9304.500599999999
$ws.Range("M122").Value = -1562.80003
$ws.Range("N122").Value = -14204.5006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2299.524
$ws.Range("I97").Value = 1986.875
$ws.Range("K97").Value = 1986.875
$ws.Range("M97").Value = -1490.875
$ws.Range("H126").Value = 2897.6
$ws.Range("I126").Value = 1579.75
$ws.Range("J126").Value = 4874.375
$ws.Range("K126").Value = 4739.25
$ws.Range("L126").Value = 14623.125
$ws.Range("M126").Value = -2269.25
$ws.Range("N126").Value = -19563.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2674.111
$ws.Range("I7").Value = 1978.2
$ws.Range("J7").Value = 3544
$ws.Range("K7").Value = 1978.2
$ws.Range("L7").Value = 3544
$ws.Range("M7").Value = -1866.2
$ws.Range("N7").Value = -3768
$ws.Range("H22").Value = 654.5454999999999
$ws.Range("I22").Value = 812
$ws.Range("J22").Value = 564.5714
$ws.Range("K22").Value = 812
$ws.Range("L22").Value = 564.5714
$ws.Range("M22").Value = -517
$ws.Range("N22").Value = -1154.5714
$ws.Range("H27").Value = 654.5454999999999
$ws.Range("I27").Value = 812
$ws.Range("J27").Value = 564.5714
$ws.Range("K27").Value = 812
$ws.Range("L27").Value = 564.5714
$ws.Range("M27").Value = -705
$ws.Range("N27").Value = -778.5714
$ws.Range("H99").Value = 52950
$ws.Range("J99").Value = 52950
$ws.Range("L99").Value = 52950
$ws.Range("N99").Value = -58940
$ws.Range("H100").Value = 4849.8335
$ws.Range("I100").Value = 2259.6
$ws.Range("J100").Value = 6700
$ws.Range("K100").Value = 2259.6
$ws.Range("L100").Value = 6700
$ws.Range("M100").Value = -1718.6
$ws.Range("N100").Value = -7782
$ws.Range("H126").Value = 2674.111
$ws.Range("I126").Value = 1978.2
$ws.Range("J126").Value = 3544
$ws.Range("K126").Value = 5934.6
$ws.Range("L126").Value = 10632
$ws.Range("M126").Value = -3464.6
$ws.Range("N126").Value = -15572

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 45031
$ws.Range("J49").Value = 45031
$ws.Range("L49").Value = 45031
$ws.Range("N49").Value = -45491
$ws.Range("H62").Value = 3542.318
$ws.Range("J62").Value = 3418.75
$ws.Range("L62").Value = 3418.75
$ws.Range("N62").Value = -4666.75
$ws.Range("H65").Value = 3542.318
$ws.Range("J65").Value = 3418.75
$ws.Range("L65").Value = 17093.75
$ws.Range("N65").Value = -23333.75
$ws.Range("H135").Value = 166703570
$ws.Range("J135").Value = 166703570
$ws.Range("L135").Value = 166703570
$ws.Range("N135").Value = -166713710

Write-Host "Applied scheduled price/profit updates across all job sheets."
